# Apply the "MarchingMasters-1013" Week's Tasks slide edits.
# Target shape: Slide 5, "Content Placeholder 2" (3rd shape on the slide).
#
# Strategy: use TextRange.Characters(start, length) sub-ranges addressed by the
# ORIGINAL (pre-edit) 1-based character offsets of the text frame, and apply
# every edit from the END of the text frame toward the START (and, within a
# paragraph that has several runs, from its LAST run toward its FIRST run).
# That way, each edit's offset is computed against text that has not yet been
# shifted by any later (numerically higher-offset) edit, so we never need to
# recompute offsets mid-script.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 12: "Goal: Research possible technologies..." -> "Goal: Prepare plan/design of P.O.C." ---
$c = $tr.Characters(637, 78)
$c.Text = "Goal: Prepare plan/design of P.O.C."

# --- Paragraph 11: "Prepare a plan on how to approach proving the concept, including proposed technologies, tasks, etc. " -> removed entirely ---
$c = $tr.Characters(536, 101)
$c.Delete()

# --- Paragraph 10: "Proof of Concept (Aparna, Brandin)" -> removed entirely ---
$c = $tr.Characters(501, 35)
$c.Delete()

# --- Paragraph 9: "Goal: Compose document or slide deck to present findings." -> removed entirely ---
$c = $tr.Characters(443, 58)
$c.Delete()

# --- Paragraph 8: "Obtain market research statistics regarding potential customers within the market."
#     -> "Prepare a plan on how to approach proving the " + "concept,based" + " on proposed technologies, tasks, etc. "
$c = $tr.Characters(360, 82)
$c.Text = "Prepare a plan on how to approach proving the concept,based on proposed technologies, tasks, etc. "
# Re-stamp the middle word so it becomes its own run (splits the single run into three).
$c = $tr.Characters(360 + 46, 13)
$c.Text = "concept,based"

# --- Paragraph 7: "Market Research (Adam, " / "Tumaris" / ")" -> "Proof of Concept (" / "Siddharth,Adam" / ")" ---
# (right-to-left within the paragraph so the first run's offset stays valid)
$c = $tr.Characters(351, 7)
$c.Text = "Siddharth,Adam"
$c = $tr.Characters(328, 23)
$c.Text = "Proof of Concept ("

# --- Paragraph 6: "Goal: Identify at least 5 users and have interviews scheduled." -> "Goal: Conduct Interviews with potential users and continue work on document." ---
$c = $tr.Characters(265, 62)
$c.Text = "Goal: Conduct Interviews with potential users and continue work on document."

# --- Paragraph 5: "Identify users of the system and Conduct interviews to determine user roles and needs. " -> removed entirely ---
$c = $tr.Characters(177, 88)
$c.Delete()

# --- Paragraph 4: "Customer Requirements (Jeffer, Brandin)" -> "Customer Requirements (Brandin,Jeffer)" ---
# (right-to-left within the paragraph so the earlier run's offset stays valid)
$c = $tr.Characters(166, 10)
$c.Text = ")"
$c = $tr.Characters(160, 6)
$c.Text = "Brandin,Jeffer"

# --- Paragraph 2: "Develop project schedule to be used for accomplishing project tasks." -> "Develop project schedule (baselined) to be used for accomplishing project tasks." ---
$c = $tr.Characters(38, 68)
$c.Text = "Develop project schedule (baselined) to be used for accomplishing project tasks."

# --- Paragraph 1: "Project Schedule (Aparna, Siddharth)" -> "Project Schedule (" + "Aparna,Tumaris" + ")" ---
$c = $tr.Characters(1, 36)
$c.Text = "Project Schedule (Aparna,Tumaris)"
# Re-stamp "Aparna,Tumaris" so it becomes its own run (splits into three runs).
$c = $tr.Characters(1 + 18, 14)
$c.Text = "Aparna,Tumaris"

# --- Turn off the autofit line-spacing reduction (normAutofit lnSpcReduction="10000" -> normAutofit) ---
$tf.AutoSize = 2
